$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextCell 'D2' '34.037.69'
Set-TextCell 'E2' '  -1.45%  '
Set-TextCell 'D3' '1.792.69'
Set-TextCell 'E3' '  -2.10%  '
Set-TextCell 'E4' '  +0.41%  '
Set-TextCell 'D5' '228.42'
Set-TextCell 'E5' '  -2.59%  '
Set-TextCell 'D6' '0.554'
Set-TextCell 'E6' '  +0.18%  '
Set-TextCell 'E7' '  +0.16%  '
Set-TextCell 'D8' '31.23'
Set-TextCell 'E8' '  -2.29%  '
Set-TextCell 'D9' '46.04'
Set-TextCell 'E9' '  -0.37%  '
Set-TextCell 'E10' '  -2.22%  '
Set-TextCell 'D11' '0.0661'
Set-TextCell 'E11' '  -3.68%  '
Set-TextCell 'D12' '0.0926'
Set-TextCell 'E12' '  -0.43%  '
Set-TextCell 'D13' '2.053.09'
Set-TextCell 'E13' '  -1.93%  '
Set-TextCell 'D14' '11.22'
Set-TextCell 'E14' '  +8.12%  '
Set-TextCell 'D15' '1.795.02'
Set-TextCell 'E15' '  -1.86%  '
Set-TextCell 'D16' '0.635'
Set-TextCell 'E16' '  -2.12%  '
Set-TextCell 'D17' '34.064.18'
Set-TextCell 'E17' '  -1.27%  '
Set-TextCell 'D18' '4.21'
Set-TextCell 'E18' '  -3.22%  '
Set-TextCell 'D19' '69.59'
Set-TextCell 'E19' '  -3.30%  '
Set-TextCell 'D20' '253.11'
Set-TextCell 'E20' '  -4.27%  '
Set-TextCell 'D21' '0.0₃0745'
Set-TextCell 'E21' '  -2.48%  '
Set-TextCell 'E22' '  +0.33%  '
Set-TextCell 'D23' '10.44'
Set-TextCell 'E23' '  -1.81%  '
Set-TextCell 'D24' '4.28'
Set-TextCell 'E24' '  -3.83%  '
Set-TextCell 'E25' '  -2.23%  '
Set-TextCell 'D26' '157.77'
Set-TextCell 'E26' '  -3.00%  '
Set-TextCell 'D27' '16.61'
Set-TextCell 'E27' '  -2.60%  '
Set-TextCell 'B28' 'Stellar'
Set-TextCell 'C28' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 'D28' '0.114'
Set-TextCell 'E28' '  -3.18%  '
Set-TextCell 'B29' 'Cosmos'
Set-TextCell 'C29' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell 'D29' '7.02'
Set-TextCell 'E29' '  -2.65%  '
Set-TextCell 'E30' '  +0.37%  '
Set-TextCell 'D31' '3.89'
Set-TextCell 'E31' '  +0.11%  '
Set-TextCell 'D32' '0.0517'
Set-TextCell 'E32' '  -0.53%  '
Set-TextCell 'E33' '  -1.04%  '
Set-TextCell 'D34' '3.64'
Set-TextCell 'E34' '  +0.68%  '
Set-TextCell 'D35' '1.84'
Set-TextCell 'E35' '  -1.39%  '
Set-TextCell 'D36' '1.476.00'
Set-TextCell 'E36' '  -8.22%  '
Set-TextCell 'E37' '  +0.22%  '
Set-TextCell 'D38' '0.635'
Set-TextCell 'E38' '  +0.03%  '
Set-TextCell 'D39' '0.0187'
Set-TextCell 'E39' '  -1.59%  '
Set-TextCell 'D40' '83.83'
Set-TextCell 'E40' '  -5.22%  '
Set-TextCell 'E41' '  +0.15%  '
Set-TextCell 'E42' '  -0.62%  '
Set-TextCell 'D43' '0.904'
Set-TextCell 'E43' '  -2.88%  '
Set-TextCell 'D44' '2.06'
Set-TextCell 'E44' '  -4.85%  '
Set-TextCell 'E45' '  -1.51%  '
Set-TextCell 'E46' '  +0.87%  '
Set-TextCell 'D47' '1.950.68'
Set-TextCell 'E47' '  -1.69%  '
Set-TextCell 'B48' 'PaxDollar'
Set-TextCell 'C48' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextCell 'D48' '1.00'
Set-TextCell 'E48' '  +0.36%  '
Set-TextCell 'B49' 'FraxShare'
Set-TextCell 'C49' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 'D49' '5.70'
Set-TextCell 'E49' '  -2.17%  '
Set-TextCell 'D50' '11.81'
Set-TextCell 'E50' '  +2.73%  '
Set-TextCell 'D51' '51.49'
Set-TextCell 'E51' '  -5.62%  '
